$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new distance column "d=6" needs to be inserted between the existing
# "d=5" (column F) and "d=7" (old column G) columns. Inserting a whole
# column at G shifts the old G (d=7) and H (d=10) data/headers one
# column to the right (G->H, H->I) and creates a fresh, empty column G
# in their place, matching the structural change in the diff.
$ws.Columns("G").Insert()

# Header for the newly inserted column
$ws.Range("G1").Value = "d=6"

# New "d=6" values for each distribution row
$ws.Range("G2").Value = 97.80141191575052
$ws.Range("G3").Value = 98.00978560161275
$ws.Range("G4").Value = 97.99836373005171
$ws.Range("G5").Value = 97.91808310958237
$ws.Range("G6").Value = 98.0118366189915
